# NC92Soil Versione 0.7 bis del 23/09/2020
# - Implementato variatore di curve di decadimento secondo il modello di Darendeli
# - Implementata possibilita' di specificare la variabilita' delle curve di
#   decadimento anche in fase di analisi batch con profili
#
# On the "Stochastic" sheet, the "Degradation curve Std" column (H) used to
# hold a plain numeric placeholder (1) for every row; it now references the
# new "Darendeli" degradation-curve-variability model (a text choice, like
# the "Inter-layer correlation" column). Row 3's "Inter-layer correlation"
# value (I3) is fixed up from a stray formula string to the correct
# "Toro: USGS AB" choice used by the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stochastic")

# Degradation curve Std (column H) now selects the Darendeli variability model
$ws.Range("H2").Value = "Darendeli"
$ws.Range("H3").Value = "Darendeli"
$ws.Range("H4").Value = "Darendeli"

# Correct the Inter-layer correlation value for row 3 to match rows 2 and 4
$ws.Range("I3").Value = "Toro: USGS AB"

# Move/restore the active selection to H5, as left by the author on save
$ws.Activate()
[void]$ws.Range("H5").Select()
